# Fruta / hortaliza, semanal
# Insert two new weekly records for "Albahaca" (date serial 44610) right
# before the existing row 372, pushing the remaining rows (old 372-398)
# down to 374-400 and extending the sheet to A1:R400.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 372 (inherits formatting, e.g. the
# date style on column D, from the surrounding rows).
$ws.Rows.Item(372).Insert()
$ws.Rows.Item(372).Insert()

$newRows = @(
    @{Row=372; I="Primera"; J=990; K=2500; L=3000; M=2677; P=446},
    @{Row=373; I="Segunda"; J=420;  K=2000; L=2000; M=2000; P=333}
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value  = 6
    $ws.Cells.Item($r, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value  = "Metropolitana"
    $ws.Cells.Item($r, 4).Value  = 44610
    $ws.Cells.Item($r, 5).Value  = 13
    $ws.Cells.Item($r, 6).Value  = 100112052
    $ws.Cells.Item($r, 7).Value  = "Albahaca"
    $ws.Cells.Item($r, 8).Value  = "Sin especificar"
    $ws.Cells.Item($r, 9).Value  = $rowData.I
    $ws.Cells.Item($r, 10).Value = $rowData.J
    $ws.Cells.Item($r, 11).Value = $rowData.K
    $ws.Cells.Item($r, 12).Value = $rowData.L
    $ws.Cells.Item($r, 13).Value = $rowData.M
    $ws.Cells.Item($r, 14).Value = "`$/docena de matas"
    $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($r, 16).Value = $rowData.P
    $ws.Cells.Item($r, 17).Value = 6
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
